# Generate Report for Handback
# Refresh the handback-status report: the old handoff files
# (133ac459-bf65-4622-bfe2-9af78ebb57f6.md / 82b55730-5605-47ed-bf66-fe9ecdfb4707.md)
# are replaced by a newer run (d29aeb6d-c59e-4536-9793-bcd535208054.md /
# ffff1c900cbd-ee74-48de-b092-3167c93c63ac.md), along with refreshed xliff
# correspondence file names and updated generate/handback timestamps.

$wb = $excel.ActiveWorkbook

# New handoff file names replacing the previous CI run's files.
$newFile1 = "d29aeb6d-c59e-4536-9793-bcd535208054.md"
$newFile2 = "ffff1c900cbd-ee74-48de-b092-3167c93c63ac.md"

# New xliff correspondence file names (per target language).
$newXlfZh = "d29aeb6d-c59e-4536-9793-bcd535208054.5afe5d86e2ef23868f0970ede3cfccfcf7ee92ef.zh-cn.xlf"
$newXlfDe = "d29aeb6d-c59e-4536-9793-bcd535208054.5afe5d86e2ef23868f0970ede3cfccfcf7ee92ef.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = $newFile1
$ovw.Range("B2").Value = "e2e\" + $newFile1
$ovw.Range("G2").Value = "2016-09-07 05:21:53"

$ovw.Range("A3").Value = $newFile2
$ovw.Range("B3").Value = "e2e\" + $newFile2
$ovw.Range("G3").Value = "2016-09-07 05:21:53"

foreach ($hl in $ovw.Hyperlinks) {
  $addr = $hl.Address
  if ($addr -like "*133ac459-bf65-4622-bfe2-9af78ebb57f6*") {
    $hl.TextToDisplay = "e2e\" + $newFile1
  } elseif ($addr -like "*82b55730-5605-47ed-bf66-fe9ecdfb4707*") {
    $hl.TextToDisplay = "e2e\" + $newFile2
  }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newFile1
$zh.Range("G2").Value = $newXlfZh
$zh.Range("H2").Value = "2016-09-07 05:21:46"
$zh.Range("I2").Value = $newFile1
$zh.Range("J2").Value = $newXlfZh
$zh.Range("K2").Value = "2016-09-07 05:22:11"

$zh.Range("A3").Value = $newFile2
$zh.Range("G3").Value = $newXlfZh
$zh.Range("H3").Value = "2016-09-07 05:21:46"
$zh.Range("I3").Value = $newFile2
$zh.Range("J3").Value = $newXlfZh
$zh.Range("K3").Value = "2016-09-07 05:22:11"

foreach ($hl in $zh.Hyperlinks) {
  $addr = $hl.Address
  if ($addr -like "*133ac459-bf65-4622-bfe2-9af78ebb57f6*") {
    $hl.TextToDisplay = $newFile1
  } elseif ($addr -like "*82b55730-5605-47ed-bf66-fe9ecdfb4707*") {
    $hl.TextToDisplay = $newFile2
  }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newFile1
$de.Range("G2").Value = $newXlfDe
$de.Range("H2").Value = "2016-09-07 05:21:53"
$de.Range("I2").Value = $newFile1
$de.Range("J2").Value = $newXlfDe
$de.Range("K2").Value = "2016-09-07 05:22:20"

$de.Range("A3").Value = $newFile2
$de.Range("G3").Value = $newXlfDe
$de.Range("H3").Value = "2016-09-07 05:21:53"
$de.Range("I3").Value = $newFile2
$de.Range("J3").Value = $newXlfDe
$de.Range("K3").Value = "2016-09-07 05:22:20"

foreach ($hl in $de.Hyperlinks) {
  $addr = $hl.Address
  if ($addr -like "*133ac459-bf65-4622-bfe2-9af78ebb57f6*") {
    $hl.TextToDisplay = $newFile1
  } elseif ($addr -like "*82b55730-5605-47ed-bf66-fe9ecdfb4707*") {
    $hl.TextToDisplay = $newFile2
  }
}
